# ---------------------------------------------------------------------------
# Applies the commit's changes to test_report.xlsx:
#   1. Update ad-copy text (column C) on "Image Alt Attribute Test".
#   2. Replace the stale menu_id in every redirect URL on "URL Status Test".
#   3. Replace the single placeholder row on "Currency Filter Test" with the
#      full set of validated currency rows.
#   4. Insert a brand-new "Script Data Extraction Test" sheet (right after
#      "Currency Filter Test").
#   5. Move "H1 Tag Existence" so it becomes the last sheet tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. "Image Alt Attribute Test" -> refreshed Spanish ad copy -----------
$imgSheet = $wb.Worksheets.Item("Image Alt Attribute Test")
$imgSheet.Range("C63").Value = "Desbloquee el mejor precio en Madrid"
$imgSheet.Range("C65").Value = "Gran oferta para esta noche"
$imgSheet.Range("C67").Value = "Descuento exclusivo para Madrid"
$imgSheet.Range("C69").Value = "Tarifas calientes hoy.  Madrid Descuentos"
$imgSheet.Range("C71").Value = "Bueno para las familias.  Libro Madrid"
$imgSheet.Range("C73").Value = "Villa de lujo a partir de Hotala™ 55"

# --- 2. "URL Status Test" -> bump the stale menu_id campaign value --------
$urlSheet = $wb.Worksheets.Item("URL Status Test")
$urlSheet.Cells.Replace("1733832290831", "1733890987713")

# --- 3. "Currency Filter Test" -> real currency validation rows -----------
$currencySheet = $wb.Worksheets.Item("Currency Filter Test")
$currencyRows = @(
    @("$ US", "Pass", "Currency $ US validated successfully."),
    @("$ CA", "Pass", "Currency $ CA validated successfully."),
    @("€ BE", "Pass", "Currency € BE validated successfully."),
    @("£ IE", "Pass", "Currency £ IE validated successfully."),
    @("$ AU", "Pass", "Currency $ AU validated successfully."),
    @("$ SG", "Pass", "Currency $ SG validated successfully."),
    @("د.إ. AE", "Pass", "Currency د.إ. AE validated successfully."),
    @("৳ BD", "Pass", "Currency ৳ BD validated successfully.")
)
$rowIndex = 2
foreach ($row in $currencyRows) {
    $currencySheet.Cells.Item($rowIndex, 1).Value = $row[0]
    $currencySheet.Cells.Item($rowIndex, 2).Value = $row[1]
    $currencySheet.Cells.Item($rowIndex, 3).Value = $row[2]
    $rowIndex = $rowIndex + 1
}

# --- 4. New sheet: "Script Data Extraction Test" ---------------------------
$scriptSheet = $wb.Worksheets.Add($null, $currencySheet)
$scriptSheet.Name = "Script Data Extraction Test"

$scriptSheet.Cells.Item(1, 1).Value = "SiteURL"
$scriptSheet.Cells.Item(1, 2).Value = "CampaignID"
$scriptSheet.Cells.Item(1, 3).Value = "SiteName"
$scriptSheet.Cells.Item(1, 4).Value = "Browser"
$scriptSheet.Cells.Item(1, 5).Value = "Country"
$scriptSheet.Cells.Item(1, 6).Value = "IP"

$scriptSheet.Cells.Item(2, 1).Value = "https://www.alojamiento.io"
$scriptSheet.Cells.Item(2, 2).Value = "ALOJAMIENTO"
$scriptSheet.Cells.Item(2, 3).Value = "alo"
$scriptSheet.Cells.Item(2, 4).Value = "Firefox"
$scriptSheet.Cells.Item(2, 5).Value = "BD"
$scriptSheet.Cells.Item(2, 6).Value = "182.160.106.203"

$scriptSheet.Range("A1:F1").Font.Bold = $true

# --- 5. Move "H1 Tag Existence" to the end of the tab strip ---------------
$h1Sheet = $wb.Worksheets.Item("H1 Tag Existence")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$h1Sheet.Move($null, $lastSheet)
